# Update countries & provincias Spain
# Applies the 8-May-2020 18:04 data refresh to the "Pais" worksheet:
#  - Updates the "last updated" timestamp in A1
#  - Refreshes case totals for several countries
#  - Two pairs of countries swapped rank (their rows exchange places)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Mayo de 2020 a las 18:04"

# --- Row 4: Estados Unidos (simple refresh, no rank change) -----------
$ws.Cells.Item(4, 2).Value = 1297003
$ws.Cells.Item(4, 3).Value = 4380
$ws.Cells.Item(4, 4).Value = 218928
$ws.Cells.Item(4, 5).Value = 1000885
$ws.Cells.Item(4, 6).Value = 16983
$ws.Cells.Item(4, 7).Value = 262
$ws.Cells.Item(4, 8).Value = 77190

# --- Rows 26/27: Chile overtakes Suecia --------------------------------
# Row 26 becomes Chile with refreshed totals
$ws.Cells.Item(26, 1).Value = "Chile"
$ws.Cells.Item(26, 2).Value = 25972
$ws.Cells.Item(26, 3).Value = 1391
$ws.Cells.Item(26, 4).Value = 12160
$ws.Cells.Item(26, 5).Value = 13518
$ws.Cells.Item(26, 6).Value = 429
$ws.Cells.Item(26, 7).Value = 9
$ws.Cells.Item(26, 8).Value = 294

# Row 27 becomes Suecia, keeping its previous (unrefreshed) totals
$ws.Cells.Item(27, 1).Value = "Suecia"
$ws.Cells.Item(27, 2).Value = 25265
$ws.Cells.Item(27, 3).Value = 642
$ws.Cells.Item(27, 4).Value = 4971
$ws.Cells.Item(27, 5).Value = 17119
$ws.Cells.Item(27, 6).Value = 425
$ws.Cells.Item(27, 7).Value = 135
$ws.Cells.Item(27, 8).Value = 3175

# --- Row 39: Banglades (simple refresh, no rank change) ----------------
$ws.Cells.Item(39, 4).Value = 2101
$ws.Cells.Item(39, 5).Value = 10827

# --- Rows 48/49: Chequia overtakes Noruega -----------------------------
# Row 48 becomes Chequia with refreshed totals
$ws.Cells.Item(48, 1).Value = "Chequia"
$ws.Cells.Item(48, 2).Value = 8065
$ws.Cells.Item(48, 3).Value = 34
$ws.Cells.Item(48, 4).Value = 4408
$ws.Cells.Item(48, 5).Value = 3386
$ws.Cells.Item(48, 6).Value = 258
$ws.Cells.Item(48, 7).Value = 1
$ws.Cells.Item(48, 8).Value = 271

# Row 49 becomes Noruega, keeping its previous (unrefreshed) totals
$ws.Cells.Item(49, 1).Value = "Noruega"
$ws.Cells.Item(49, 2).Value = 8055
$ws.Cells.Item(49, 3).Value = 21
$ws.Cells.Item(49, 4).Value = 32
$ws.Cells.Item(49, 5).Value = 7805
$ws.Cells.Item(49, 6).Value = 27
$ws.Cells.Item(49, 7).Value = 1
$ws.Cells.Item(49, 8).Value = 218

# --- Row 58: Argelia (simple refresh, no rank change) -------------------
$ws.Cells.Item(58, 2).Value = 5369
$ws.Cells.Item(58, 3).Value = 187
$ws.Cells.Item(58, 4).Value = 2467
$ws.Cells.Item(58, 5).Value = 2414
$ws.Cells.Item(58, 7).Value = 5
$ws.Cells.Item(58, 8).Value = 488

# --- Row 70: Grecia (simple refresh, no rank change) ---------------------
$ws.Cells.Item(70, 2).Value = 2691
$ws.Cells.Item(70, 3).Value = 13
$ws.Cells.Item(70, 5).Value = 1167
$ws.Cells.Item(70, 6).Value = 32
$ws.Cells.Item(70, 7).Value = 2
$ws.Cells.Item(70, 8).Value = 150

# --- Rows 192/193: Nueva Caledonia overtakes Belice (totals unchanged) --
$ws.Cells.Item(192, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(192, 4).Value = 18
$ws.Cells.Item(192, 8).Value = 0

$ws.Cells.Item(193, 1).Value = "Belice"
$ws.Cells.Item(193, 4).Value = 16
$ws.Cells.Item(193, 8).Value = 2
